$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Update the "Section N" header cells in column A of the first sheet.
# (mirrors the bug fix described in the commit message, where the section
# numbers parsed by XlsFileParserService changed, and one of them now
# resolves to a plain number 0 instead of a string)
$ws1.Range("A1").Value = "Section 1123"
$ws1.Range("A3").Value = "Section 243"
$ws1.Range("A6").Value = "Section 365"
$ws1.Range("A12").Value = "Section 44"
$ws1.Range("A15").Value = "Section 544"
$ws1.Range("A18").Value = 0
$ws1.Range("A21").Value = "Section 7908"
$ws1.Range("A26").Value = "Section 889087"

# Update the active sheet / selection state on each sheet so that the
# first worksheet becomes the active / selected tab (instead of the
# second one) and the cell selections match the new view state.
$ws2.Activate() | Out-Null
$ws2.Range("C5").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("J9").Select() | Out-Null
